$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "CALL /path/1"
$ws.Range("G3").Value = "CALL /path/2"
$ws.Range("G4").Value = "CALL /path/3"

$ws.Range("G4").Select()
